# Add two new paragraphs ("Login" and a tab-indented description) right
# after the existing "Driver for program" paragraph, moving the
# _GoBack bookmark so it ends up trailing the run in the new last
# paragraph (matching the target OOXML).

$d = $word.ActiveDocument

# The document's last paragraph ("Driver for program") currently owns the
# trailing _GoBack bookmark. Its (collapsed) paragraph-mark range sits at
# the very end of the story; inserting the new paragraphs' text in one
# shot immediately before that mark pushes the bookmark onto the newly
# created last paragraph, while leaving "Driver for program" itself
# untouched.
$endPos = $d.Content.End
$endMark = $d.Range($endPos - 1, $endPos - 1)
$endMark.InsertBefore("`rLogin`r")

# Fill in the final paragraph's text (tab character + description) by
# replacing the whole paragraph's range, which keeps the bookmark
# anchored after the run, same as in the target document.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = [char]9 + "Checks for user info and logs them into their account"
